$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TPM-derived values on row 2 (Pomc-Mc2r pair) per new script output
$ws.Range("G2").Value = 1.524170333333333
$ws.Range("H2").Value = 4.572511

$ws.Range("Q2").Value = 0.051861419762
$ws.Range("R2").Value = 0.4667527778580001
